$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.776.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "'2.475.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'576.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'149.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.543"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").Value = "'2.474.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "'5.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "'0.360"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "'27.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "'0.0000182"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "'2.921.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "'63.620.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "'2.464.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").Value = "'11.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").Value = "'7.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.78%  "
$ws.Range("D21").Value = "'330.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "'2.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.86%  "
$ws.Range("D24").Value = "'0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'65.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D26").Value = "'628.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.17%  "
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("D28").Value = "'8.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "'2.599.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.87%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'8.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").Value = "'5.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.22%  "
$ws.Range("D36").Value = "'1.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'0.385"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "'5.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").Value = "'18.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.09%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'147.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'1.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "'151.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").Value = "'21.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").Value = "'0.606"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").Value = "'0.0922"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.69%  "
